$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.018.01'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = '1.642.83'
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0639'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.869.69'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = '1.675.20'
$ws.Range("E14").Value = '  +2.81%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.02'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '25.930.55'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("E25").Value = '  +5.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.70'
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.92'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("E35").Value = '  +2.54%  '
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").Value = '1.134.44'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.51'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.61'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.798'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '1.778.95'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E45").Value = '  +3.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.79'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0530'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.71%  '
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("E51").Value = '  -0.21%  '
